$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial number that was bumped by one
# day (46060 -> 46061) for every data row (rows 2 through 99).
$ws.Range("C2:C99").Value = 46061
